$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the small block of rows 6-11 (columns A:B) by column A, ascending,
# matching the re-ordering performed in Excel by the lab member.
$sortRange = $ws.Range("A6:B11")
$keyRange = $ws.Range("A6:A11")
[void]$sortRange.Sort($keyRange, 1)

# Update the active selection to reflect the cell last interacted with.
[void]$ws.Range("J1:J22").Select()
